$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "(N, 0)" placeholder strings with "(N, None)" in column A
$cells = @("A6","A7","A11","A13","A14","A16","A18","A19","A22","A24","A25","A28","A29","A30","A32","A34")
foreach ($addr in $cells) {
    $cell = $ws.Range($addr)
    $old = $cell.Value2
    $new = $old -replace ", 0\)$", ", None)"
    $cell.Value = $new
}

# Increase the height of row 4 to accommodate the new comment/content
$ws.Rows.Item(4).RowHeight = 95.25
